$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

# Update K2 (video_Interior for row 2)
$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109864/1750183612496/A5MTrC4t3W-camera-video-segment-driver-1750183614996.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBYDP2TS4%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T180004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKr%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIFpWzNZHTIuWx9Jfzil6EW78cqt8sJ%2FyIKXs7AqYf58IAiEA%2FplU9EnLA%2FBKbqOQgPr%2BpbyuqGuuuXTTshUMwZx9cVMq5gMIk%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDGPA%2FyuA9QEJpPr9DCq6Azw6oOyCKDKWrzL%2BHT%2FvK6UaJ2Lzu5Cf4R24oyE%2FlinXD0XsbTuWcmyYsn33d0ziGODpr1IxKv2fGo8pxPMoRS88PhZQOKtx3bXwg9nww1G8NH5neYXf16JZbDbTVVoU%2FAO5jWTR%2BR0%2FscWsGnNq%2FDPDnG%2FDp4ZD4o47mwcRExrnTCRFKs3Bjoi9ZRlko16yUzNitA8Pc8TPl3XgutBSpzQeAnADoQC44bB512exlSR8nRx6vJwvniQv3EAizZwtO8xT8oW%2BU3GkpXLOEa1Ql7u0HXNTl7WRSUm0ZqsLkGHKwxkwkPNy9fCheupmAAewRmeUMF6J4KnidkaH8WdLTp9urm1cFN4Kuutl2npiaOS%2FnhxV%2F3A4SN7XjRZKX89F3uytxQnl%2BjFW9%2BiK%2Ba6tBb34%2FM95IEXlATDB1rHKO5rnuhzqk784KOb7KXNUs76pLZasYyBFfU6s7wVXlO8VQdVI8I56j9JDuOFwqD9Vc%2BpRCp61IkkIp%2BQVIrS0l9VyiFge3JRBXfhU%2BwW94c2VOOSZkCGKotw%2Bp8uzJ9uaXeGYipFxMoLH3JHNgF1BKE42l2TWGb8gjLZBFuswtffLwgY6pQFdUfnEVvfK7XJ%2FoXg%2FIdPc5K9ThBH%2BcgN94Xt8N3f%2BIVTYdzorSig1ONBT2GLTJbtbw84PvW1lzMJ1EDChZEoUXlrDo70l%2FvIG83IdRKzCx0MDkGTDAv4d8RCsbOAEDU9HZspWl%2B9ePS%2BqJsMUCjAR6YmnsoU7AQAa0BpKJOP8M2K2jIJ0VIyw5K2siV7bkhkj%2BhunwNDvx3pfcOTLVtPxL6wEh8E%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2002%3A00%3A04%20GMT&X-Amz-Signature=81c806f74217044d5203a46cb5122592e381dac49a440e6a9cf9404fc6f74881'

# Update L2 (video_Exterior for row 2)
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1750183612496/TAdN8NokD3-camera-video-segment-1750183614996.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBYDP2TS4%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T180004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKr%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIFpWzNZHTIuWx9Jfzil6EW78cqt8sJ%2FyIKXs7AqYf58IAiEA%2FplU9EnLA%2FBKbqOQgPr%2BpbyuqGuuuXTTshUMwZx9cVMq5gMIk%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDGPA%2FyuA9QEJpPr9DCq6Azw6oOyCKDKWrzL%2BHT%2FvK6UaJ2Lzu5Cf4R24oyE%2FlinXD0XsbTuWcmyYsn33d0ziGODpr1IxKv2fGo8pxPMoRS88PhZQOKtx3bXwg9nww1G8NH5neYXf16JZbDbTVVoU%2FAO5jWTR%2BR0%2FscWsGnNq%2FDPDnG%2FDp4ZD4o47mwcRExrnTCRFKs3Bjoi9ZRlko16yUzNitA8Pc8TPl3XgutBSpzQeAnADoQC44bB512exlSR8nRx6vJwvniQv3EAizZwtO8xT8oW%2BU3GkpXLOEa1Ql7u0HXNTl7WRSUm0ZqsLkGHKwxkwkPNy9fCheupmAAewRmeUMF6J4KnidkaH8WdLTp9urm1cFN4Kuutl2npiaOS%2FnhxV%2F3A4SN7XjRZKX89F3uytxQnl%2BjFW9%2BiK%2Ba6tBb34%2FM95IEXlATDB1rHKO5rnuhzqk784KOb7KXNUs76pLZasYyBFfU6s7wVXlO8VQdVI8I56j9JDuOFwqD9Vc%2BpRCp61IkkIp%2BQVIrS0l9VyiFge3JRBXfhU%2BwW94c2VOOSZkCGKotw%2Bp8uzJ9uaXeGYipFxMoLH3JHNgF1BKE42l2TWGb8gjLZBFuswtffLwgY6pQFdUfnEVvfK7XJ%2FoXg%2FIdPc5K9ThBH%2BcgN94Xt8N3f%2BIVTYdzorSig1ONBT2GLTJbtbw84PvW1lzMJ1EDChZEoUXlrDo70l%2FvIG83IdRKzCx0MDkGTDAv4d8RCsbOAEDU9HZspWl%2B9ePS%2BqJsMUCjAR6YmnsoU7AQAa0BpKJOP8M2K2jIJ0VIyw5K2siV7bkhkj%2BhunwNDvx3pfcOTLVtPxL6wEh8E%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2002%3A00%3A04%20GMT&X-Amz-Signature=8eeb7cc11346edaffdefa1b327dfa5b90317d01851508e24d6c26bee5437fc31'

# Update K3 (video_Interior for row 3)
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991152736/1750182961683/iBjUYEdqdu-camera-video-segment-driver-1750182964183.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBYDP2TS4%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T180004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKr%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIFpWzNZHTIuWx9Jfzil6EW78cqt8sJ%2FyIKXs7AqYf58IAiEA%2FplU9EnLA%2FBKbqOQgPr%2BpbyuqGuuuXTTshUMwZx9cVMq5gMIk%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDGPA%2FyuA9QEJpPr9DCq6Azw6oOyCKDKWrzL%2BHT%2FvK6UaJ2Lzu5Cf4R24oyE%2FlinXD0XsbTuWcmyYsn33d0ziGODpr1IxKv2fGo8pxPMoRS88PhZQOKtx3bXwg9nww1G8NH5neYXf16JZbDbTVVoU%2FAO5jWTR%2BR0%2FscWsGnNq%2FDPDnG%2FDp4ZD4o47mwcRExrnTCRFKs3Bjoi9ZRlko16yUzNitA8Pc8TPl3XgutBSpzQeAnADoQC44bB512exlSR8nRx6vJwvniQv3EAizZwtO8xT8oW%2BU3GkpXLOEa1Ql7u0HXNTl7WRSUm0ZqsLkGHKwxkwkPNy9fCheupmAAewRmeUMF6J4KnidkaH8WdLTp9urm1cFN4Kuutl2npiaOS%2FnhxV%2F3A4SN7XjRZKX89F3uytxQnl%2BjFW9%2BiK%2Ba6tBb34%2FM95IEXlATDB1rHKO5rnuhzqk784KOb7KXNUs76pLZasYyBFfU6s7wVXlO8VQdVI8I56j9JDuOFwqD9Vc%2BpRCp61IkkIp%2BQVIrS0l9VyiFge3JRBXfhU%2BwW94c2VOOSZkCGKotw%2Bp8uzJ9uaXeGYipFxMoLH3JHNgF1BKE42l2TWGb8gjLZBFuswtffLwgY6pQFdUfnEVvfK7XJ%2FoXg%2FIdPc5K9ThBH%2BcgN94Xt8N3f%2BIVTYdzorSig1ONBT2GLTJbtbw84PvW1lzMJ1EDChZEoUXlrDo70l%2FvIG83IdRKzCx0MDkGTDAv4d8RCsbOAEDU9HZspWl%2B9ePS%2BqJsMUCjAR6YmnsoU7AQAa0BpKJOP8M2K2jIJ0VIyw5K2siV7bkhkj%2BhunwNDvx3pfcOTLVtPxL6wEh8E%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2002%3A00%3A04%20GMT&X-Amz-Signature=39ebdf953a73ebc9c8d6182d44193a7c7f634e4d1d9bb149a2467f48845581b4'

# Update L3 (video_Exterior for row 3)
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991152736/1750182961683/Foph7ZSROm-camera-video-segment-1750182964183.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBYDP2TS4%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T180004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKr%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIFpWzNZHTIuWx9Jfzil6EW78cqt8sJ%2FyIKXs7AqYf58IAiEA%2FplU9EnLA%2FBKbqOQgPr%2BpbyuqGuuuXTTshUMwZx9cVMq5gMIk%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDGPA%2FyuA9QEJpPr9DCq6Azw6oOyCKDKWrzL%2BHT%2FvK6UaJ2Lzu5Cf4R24oyE%2FlinXD0XsbTuWcmyYsn33d0ziGODpr1IxKv2fGo8pxPMoRS88PhZQOKtx3bXwg9nww1G8NH5neYXf16JZbDbTVVoU%2FAO5jWTR%2BR0%2FscWsGnNq%2FDPDnG%2FDp4ZD4o47mwcRExrnTCRFKs3Bjoi9ZRlko16yUzNitA8Pc8TPl3XgutBSpzQeAnADoQC44bB512exlSR8nRx6vJwvniQv3EAizZwtO8xT8oW%2BU3GkpXLOEa1Ql7u0HXNTl7WRSUm0ZqsLkGHKwxkwkPNy9fCheupmAAewRmeUMF6J4KnidkaH8WdLTp9urm1cFN4Kuutl2npiaOS%2FnhxV%2F3A4SN7XjRZKX89F3uytxQnl%2BjFW9%2BiK%2Ba6tBb34%2FM95IEXlATDB1rHKO5rnuhzqk784KOb7KXNUs76pLZasYyBFfU6s7wVXlO8VQdVI8I56j9JDuOFwqD9Vc%2BpRCp61IkkIp%2BQVIrS0l9VyiFge3JRBXfhU%2BwW94c2VOOSZkCGKotw%2Bp8uzJ9uaXeGYipFxMoLH3JHNgF1BKE42l2TWGb8gjLZBFuswtffLwgY6pQFdUfnEVvfK7XJ%2FoXg%2FIdPc5K9ThBH%2BcgN94Xt8N3f%2BIVTYdzorSig1ONBT2GLTJbtbw84PvW1lzMJ1EDChZEoUXlrDo70l%2FvIG83IdRKzCx0MDkGTDAv4d8RCsbOAEDU9HZspWl%2B9ePS%2BqJsMUCjAR6YmnsoU7AQAa0BpKJOP8M2K2jIJ0VIyw5K2siV7bkhkj%2BhunwNDvx3pfcOTLVtPxL6wEh8E%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2002%3A00%3A04%20GMT&X-Amz-Signature=7c619e2e5927d45357fcbc6c7e2a5cc37103c33e0a2af037e93eca287a256268'

